$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old row 2 (the "And" bug entry) entirely; this shifts rows 3 & 4 up to 2 & 3.
$ws.Rows.Item(2).Delete()

# Fill in the new row 4 with the "line to line" bug entry.
$ws.Range("B4").Value = "When connecting a line to a line it goes to the top left (some times Connects with an empty input)"
$ws.Range("C4").Value = "n.a."
$ws.Range("D4").Value = "n.a."
$ws.Range("E4").Value = "When connecting a line to a line it goes to the top left (some times Connects with an empty input)"
$ws.Range("F4").Value = "n.a."

$ws.Range("C28").Select()
